$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 545, shifting existing rows 545:602 down to 546:603
$ws.Rows.Item(545).Insert()

# Populate the newly inserted row 545 with the new record
$ws.Range("A545").Value = 5
$ws.Range("B545").Value = "Macroferia Regional de Talca"
$ws.Range("C545").Value = "Maule"
$ws.Range("D545").Value = 45194
$ws.Range("D545").NumberFormat = $ws.Range("D546").NumberFormat
$ws.Range("E545").Value = 7
$ws.Range("F545").Value = 100114013
$ws.Range("G545").Value = "Zanahoria"
$ws.Range("H545").Value = "Sin especificar"
$ws.Range("I545").Value = "Primera"
$ws.Range("J545").Value = 600
$ws.Range("K545").Value = 8000
$ws.Range("L545").Value = 8000
$ws.Range("M545").Value = 8000
$ws.Range("N545").Value = "$/saco 20 kilos"
$ws.Range("O545").Value = "Región de Ñuble"
$ws.Range("P545").Value = 400
$ws.Range("Q545").Value = 20
$ws.Range("R545").Value = "Hortaliza"
